# timeLog.xlsx - "Player movement - improved :D (final?)"
#
# A new time-tracking entry for Aris (17:30 - 18:40, on the existing
# 45315 date) is logged into the previously-blank row 7 of the "Tabelle1"
# sheet, and the A2:F11 block is then re-sorted by Name (column A) - the
# same "sort by name" the author had already applied to the sheet before
# (see the pre-existing <sortState> element), just re-run now that there
# is a new row to fold in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- 1. Fill in the previously empty row 7 with the new entry -------------
# Copy the formatting of an existing data row down onto row 7 first (so the
# date cell B7 picks up the date number format, matching B2:B6), then enter
# the new entry's values/formulas.
$ws.Range("A4:F4").Copy()
$ws.Range("A7:F7").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A7").Value2 = "Aris"
$ws.Range("B7").Value2 = 45315
$ws.Range("C7").Formula = "=17 + 30/60"
$ws.Range("D7").Formula = "=18 + 40/60"
$ws.Range("E7").Formula = "=D7-C7"
$ws.Range("F7").ClearContents()

# --- 2. Sort A2:F11 ascending by column A (Name), header excluded ---------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A11")) | Out-Null
$ws.Sort.SetRange($ws.Range("A1:F11"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# --- 3. The sort leaves two newly-adjacent same-name rows with a stray,
#        formatted-but-empty G cell (matches the author's sheet) ----------
$ws.Range("E74").Copy()
$ws.Range("G5").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("G6").PasteSpecial(-4122)      # xlPasteFormats

# --- 4. Re-point the saved cursor position, matching the author's sheet ---
$ws.Range("K11").Select()

Write-Host "done"
